# Auto-generated script applying meteocat daily summary refresh
# Commit: Update automàtic: dades i banners [2026-02-21 20:49]
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-21 20:48:14"
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "44%"
$ws.Range("O2").Value = "3.9 °C"
$ws.Range("E3").Value = "2026-02-21 20:48:17"
$ws.Range("K3").Value = "16.1 MJ/m2"
$ws.Range("E4").Value = "2026-02-21 20:48:19"
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "72%"
$ws.Range("O4").Value = "9.5 °C"
$ws.Range("E5").Value = "2026-02-21 20:48:21"
$ws.Range("O5").Value = "3.7 °C"
$ws.Range("E6").Value = "2026-02-21 20:48:23"
$ws.Range("E7").Value = "2026-02-21 20:48:26"
$ws.Range("E8").Value = "2026-02-21 20:48:28"
$ws.Range("K8").Value = "16.0 MJ/m2"
$ws.Range("O8").Value = "11.1 °C"
$ws.Range("E9").Value = "2026-02-21 20:48:31"
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "55%"
$ws.Range("N9").Value = "7.3 °C 20:28 TU"
$ws.Range("O9").Value = "13.5 °C"
$ws.Range("E10").Value = "2026-02-21 20:48:32"
$ws.Range("E11").Value = "2026-02-21 20:48:33"
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = "52%"
$ws.Range("O11").Value = "9.0 °C"
$ws.Range("E12").Value = "2026-02-21 20:48:34"
$ws.Range("H12").NumberFormat = "@"
$ws.Range("H12").Value = "61%"
$ws.Range("O12").Value = "12.8 °C"
$ws.Range("E13").Value = "2026-02-21 20:48:35"
$ws.Range("E14").Value = "2026-02-21 20:48:36"
$ws.Range("O14").Value = "11.4 °C"
$ws.Range("E15").Value = "2026-02-21 20:48:37"
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = "55%"
$ws.Range("O15").Value = "13.2 °C"
$ws.Range("E16").Value = "2026-02-21 20:48:38"
$ws.Range("H16").NumberFormat = "@"
$ws.Range("H16").Value = "33%"
$ws.Range("E17").Value = "2026-02-21 20:48:39"
$ws.Range("O17").Value = "8.7 °C"
$ws.Range("E18").Value = "2026-02-21 20:48:41"
$ws.Range("H18").NumberFormat = "@"
$ws.Range("H18").Value = "74%"
$ws.Range("K18").Value = "15.3 MJ/m2"
$ws.Range("O18").Value = "8.7 °C"
$ws.Range("E19").Value = "2026-02-21 20:48:42"
$ws.Range("E20").Value = "2026-02-21 20:48:43"
$ws.Range("O20").Value = "3.0 °C"
$ws.Range("E21").Value = "2026-02-21 20:48:45"
$ws.Range("E22").Value = "2026-02-21 20:48:48"
$ws.Range("E23").Value = "2026-02-21 20:48:50"
$ws.Range("E24").Value = "2026-02-21 20:48:52"
$ws.Range("H24").NumberFormat = "@"
$ws.Range("H24").Value = "83%"
$ws.Range("O24").Value = "6.4 °C"
$ws.Range("E25").Value = "2026-02-21 20:48:55"
$ws.Range("E26").Value = "2026-02-21 20:48:57"
$ws.Range("E27").Value = "2026-02-21 20:49:00"
$ws.Range("E28").Value = "2026-02-21 20:49:02"
$ws.Range("O28").Value = "8.3 °C"
$ws.Range("E29").Value = "2026-02-21 20:49:05"
$ws.Range("H29").NumberFormat = "@"
$ws.Range("H29").Value = "66%"
$ws.Range("N29").Value = "6.5 °C 20:28 TU"
$ws.Range("O29").Value = "11.7 °C"
$ws.Range("E30").Value = "2026-02-21 20:49:07"
$ws.Range("O30").Value = "11.6 °C"
$ws.Range("E31").Value = "2026-02-21 20:49:10"
$ws.Range("J31").Value = "1028.4 hPa"
$ws.Range("E32").Value = "2026-02-21 20:49:12"
$ws.Range("H32").NumberFormat = "@"
$ws.Range("H32").Value = "80%"
$ws.Range("O32").Value = "5.3 °C"
$ws.Range("E33").Value = "2026-02-21 20:49:14"
$ws.Range("E34").Value = "2026-02-21 20:49:17"
$ws.Range("H34").NumberFormat = "@"
$ws.Range("H34").Value = "39%"
$ws.Range("N34").Value = "-0.3 °C 20:17 TU"
$ws.Range("O34").Value = "4.5 °C"
$ws.Range("E35").Value = "2026-02-21 20:49:20"
$ws.Range("E36").Value = "2026-02-21 20:49:22"
$ws.Range("J36").Value = "1029.2 hPa"
$ws.Range("E37").Value = "2026-02-21 20:49:24"
$ws.Range("E38").Value = "2026-02-21 20:49:27"
$ws.Range("E39").Value = "2026-02-21 20:49:29"
$ws.Range("H39").NumberFormat = "@"
$ws.Range("H39").Value = "33%"
$ws.Range("O39").Value = "2.3 °C"
$ws.Range("E40").Value = "2026-02-21 20:49:32"
$ws.Range("H40").NumberFormat = "@"
$ws.Range("H40").Value = "52%"
$ws.Range("O40").Value = "8.8 °C"
$ws.Range("E41").Value = "2026-02-21 20:49:34"
$ws.Range("H41").NumberFormat = "@"
$ws.Range("H41").Value = "68%"
$ws.Range("E42").Value = "2026-02-21 20:49:37"
$ws.Range("H42").NumberFormat = "@"
$ws.Range("H42").Value = "74%"
$ws.Range("E43").Value = "2026-02-21 20:49:39"
$ws.Range("E44").Value = "2026-02-21 20:49:41"
$ws.Range("H44").NumberFormat = "@"
$ws.Range("H44").Value = "40%"
$ws.Range("N44").Value = "-0.8 °C 20:19 TU"
$ws.Range("O44").Value = "2.3 °C"
$ws.Range("E45").Value = "2026-02-21 20:49:44"
$ws.Range("E46").Value = "2026-02-21 20:49:46"
$ws.Range("O46").Value = "9.9 °C"
